$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.035492529207973
$ws.Range("D2").Value = 1.037378459628777
$ws.Range("E2").Value = 1.043347320323245
$ws.Range("F2").Value = 1.051433521651894
$ws.Range("I2").Value = 1.036698782795838
$ws.Range("J2").Value = 1.040606006761255
$ws.Range("K2").Value = 1.040169514315052
$ws.Range("L2").Value = 1.046121442126428
$ws.Range("M2").Value = 1.054185036440427
$ws.Range("N2").Value = 1.005712725503983
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.036345533137149
$ws.Range("D3").Value = 1.037993403367823
$ws.Range("E3").Value = 1.04414007371606
$ws.Range("F3").Value = 1.052386945588968
$ws.Range("I3").Value = 1.036872791389277
$ws.Range("J3").Value = 1.041103068719073
$ws.Range("K3").Value = 1.040594769499334
$ws.Range("L3").Value = 1.046725259811469
$ws.Range("M3").Value = 1.054950743636265
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.036898069220729
$ws.Range("D4").Value = 1.03839175007892
$ws.Range("E4").Value = 1.04465399806848
$ws.Range("F4").Value = 1.053005191315282
$ws.Range("I4").Value = 1.036984456721124
$ws.Range("J4").Value = 1.041424628648928
$ws.Range("K4").Value = 1.040869681295405
$ws.Range("L4").Value = 1.0471162729038
$ws.Range("M4").Value = 1.055446898853029
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.037130494042757
$ws.Range("D5").Value = 1.038559317984575
$ws.Range("E5").Value = 1.04487027995685
$ws.Range("F5").Value = 1.053265415280732
$ws.Range("I5").Value = 1.037031177565089
$ws.Range("J5").Value = 1.041559793909641
$ws.Range("K5").Value = 1.040985191581794
$ws.Range("L5").Value = 1.047280725910375
$ws.Range("M5").Value = 1.055655646497785
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.037169527279692
$ws.Range("D6").Value = 1.038587459361664
$ws.Range("E6").Value = 1.044906607963213
$ws.Range("F6").Value = 1.053309126360177
$ws.Range("I6").Value = 1.037039009089341
$ws.Range("J6").Value = 1.041582487629904
$ws.Range("K6").Value = 1.041004582560401
$ws.Range("L6").Value = 1.047308342423832
$ws.Range("M6").Value = 1.055690705717544
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.036901174349661
$ws.Range("D7").Value = 1.038393988727059
$ws.Range("E7").Value = 1.044656887143564
$ws.Range("F7").Value = 1.05300866721115
$ws.Range("I7").Value = 1.036985081885403
$ws.Range("J7").Value = 1.041426434808088
$ws.Range("K7").Value = 1.040871224997074
$ws.Range("L7").Value = 1.047118470055241
$ws.Range("M7").Value = 1.055449687506835
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.03578068373358
$ws.Range("D8").Value = 1.037586191259892
$ws.Range("E8").Value = 1.043615035262156
$ws.Range("F8").Value = 1.051755462050887
$ws.Range("I8").Value = 1.036757781957621
$ws.Range("J8").Value = 1.04077400544365
$ws.Range("K8").Value = 1.040313283754073
$ws.Range("L8").Value = 1.046325441652573
$ws.Range("M8").Value = 1.054443666671308
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.033810787578096
$ws.Range("D9").Value = 1.0361661694757
$ws.Range("E9").Value = 1.041786586476362
$ws.Range("F9").Value = 1.049557316380466
$ws.Range("I9").Value = 1.036350157612704
$ws.Range("J9").Value = 1.039623840282886
$ws.Range("K9").Value = 1.039328205282475
$ws.Range("L9").Value = 1.044930402497078
$ws.Range("M9").Value = 1.052676293576075
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.032500673915796
$ws.Range("D10").Value = 1.03522188756141
$ws.Range("E10").Value = 1.040572716346617
$ws.Range("F10").Value = 1.048098822853907
$ws.Range("I10").Value = 1.036073678063231
$ws.Range("J10").Value = 1.03885679811576
$ws.Range("K10").Value = 1.038670268545391
$ws.Range("L10").Value = 1.044002056308442
$ws.Range("M10").Value = 1.051501743492426
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.031934146405744
$ws.Range("D11").Value = 1.034813593157688
$ws.Range("E11").Value = 1.040048327328079
$ws.Range("F11").Value = 1.047468946416308
$ws.Range("I11").Value = 1.035952845644265
$ws.Range("J11").Value = 1.038524612513418
$ws.Range("K11").Value = 1.038385100500157
$ws.Range("L11").Value = 1.043600487656663
$ws.Range("M11").Value = 1.050994046991437
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.03172382860758
$ws.Range("D12").Value = 1.034662024161506
$ws.Range("E12").Value = 1.039853731642988
$ws.Range("F12").Value = 1.047235233566523
$ws.Range("I12").Value = 1.035907796237382
$ws.Range("J12").Value = 1.038401217491344
$ws.Range("K12").Value = 1.038279135934239
$ws.Range("L12").Value = 1.043451390352899
$ws.Range("M12").Value = 1.050805601392408
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.0317689372273
$ws.Range("D13").Value = 1.034694532138193
$ws.Range("E13").Value = 1.039895464644745
$ws.Range("F13").Value = 1.047285354365573
$ws.Range("I13").Value = 1.035917467028809
$ws.Range("J13").Value = 1.038427686417167
$ws.Range("K13").Value = 1.03830186750524
$ws.Range("L13").Value = 1.043483369349959
$ws.Range("M13").Value = 1.050846017456993
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.031916759101055
$ws.Range("D14").Value = 1.034801062577
$ws.Range("E14").Value = 1.040032238188849
$ws.Range("F14").Value = 1.047449622502426
$ws.Range("I14").Value = 1.035949125246541
$ws.Range("J14").Value = 1.038514412767995
$ws.Range("K14").Value = 1.038376342254599
$ws.Range("L14").Value = 1.043588161930384
$ws.Range("M14").Value = 1.050978467246528
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.032007852418607
$ws.Range("D15").Value = 1.034866711436508
$ws.Range("E15").Value = 1.040116533548475
$ws.Range("F15").Value = 1.047550866876499
$ws.Range("I15").Value = 1.035968608827165
$ws.Range("J15").Value = 1.038567846910004
$ws.Range("K15").Value = 1.038422223279854
$ws.Range("L15").Value = 1.043652736509828
$ws.Range("M15").Value = 1.051060091933203
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.032538288629306
$ws.Range("D16").Value = 1.035248997205655
$ws.Range("E16").Value = 1.040607544271398
$ws.Range("F16").Value = 1.048140660823017
$ws.Range("I16").Value = 1.036081673876986
$ws.Range("J16").Value = 1.038878843202602
$ws.Range("K16").Value = 1.038689188471457
$ws.Range("L16").Value = 1.044028715881339
$ws.Range("M16").Value = 1.051535456551647
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.032871221882252
$ws.Range("D17").Value = 1.035488953030921
$ws.Range("E17").Value = 1.040915871455385
$ws.Range("F17").Value = 1.048511068901719
$ws.Range("I17").Value = 1.036152298544616
$ws.Range("J17").Value = 1.039073910081683
$ws.Range("K17").Value = 1.038856575324956
$ws.Range("L17").Value = 1.044264668910469
$ws.Range("M17").Value = 1.051833879980042
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.033065489477312
$ws.Range("D18").Value = 1.035628971620669
$ws.Range("E18").Value = 1.041095831540077
$ws.Range("F18").Value = 1.048727281799222
$ws.Range("I18").Value = 1.036193385098729
$ws.Range("J18").Value = 1.039187684263061
$ws.Range("K18").Value = 1.03895418243012
$ws.Range("L18").Value = 1.044402335920105
$ws.Range("M18").Value = 1.052008031222919
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.033131742087444
$ws.Range("D19").Value = 1.035676723858841
$ws.Range("E19").Value = 1.041157213249613
$ws.Range("F19").Value = 1.048801031917026
$ws.Range("I19").Value = 1.036207376278246
$ws.Range("J19").Value = 1.039226477429063
$ws.Range("K19").Value = 1.038987459346676
$ws.Range("L19").Value = 1.044449283487929
$ws.Range("M19").Value = 1.052067426820041
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.032835493706391
$ws.Range("D20").Value = 1.035463202194086
$ws.Range("E20").Value = 1.040882778643018
$ws.Range("F20").Value = 1.048471311044364
$ws.Range("I20").Value = 1.036144732313139
$ws.Range("J20").Value = 1.039052981763113
$ws.Range("K20").Value = 1.03883861906479
$ws.Range("L20").Value = 1.044239349277611
$ws.Range("M20").Value = 1.051801853073517
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.031873226036362
$ws.Range("D21").Value = 1.034769689527489
$ws.Range("E21").Value = 1.039991956647671
$ws.Range("F21").Value = 1.047401242678897
$ws.Range("I21").Value = 1.035939807291643
$ws.Range("J21").Value = 1.038488874188984
$ws.Range("K21").Value = 1.038354412404526
$ws.Range("L21").Value = 1.043557301372442
$ws.Range("M21").Value = 1.050939460332115
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.031268880567684
$ws.Range("D22").Value = 1.034334170547912
$ws.Range("E22").Value = 1.039432937025777
$ws.Range("F22").Value = 1.046729903077205
$ws.Range("I22").Value = 1.035809997678911
$ws.Range("J22").Value = 1.038134160365393
$ws.Range("K22").Value = 1.038049738969789
$ws.Range("L22").Value = 1.043128836964105
$ws.Range("M22").Value = 1.050398024444218
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.031589191483129
$ws.Range("D23").Value = 1.034564997545942
$ws.Range("E23").Value = 1.039729181376466
$ws.Range("F23").Value = 1.047085654405642
$ws.Range("I23").Value = 1.035878903477798
$ws.Range("J23").Value = 1.038322203979167
$ws.Range("K23").Value = 1.038211273862201
$ws.Range("L23").Value = 1.043355938886696
$ws.Range("M23").Value = 1.050684974955321
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.032851637504967
$ws.Range("D24").Value = 1.035474837714719
$ws.Range("E24").Value = 1.04089773149591
$ws.Range("F24").Value = 1.048489275414177
$ws.Range("I24").Value = 1.036148151500061
$ws.Range("J24").Value = 1.039062438385108
$ws.Range("K24").Value = 1.038846732808836
$ws.Range("L24").Value = 1.044250790007637
$ws.Range("M24").Value = 1.051816324388743
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.034319504291269
$ws.Range("D25").Value = 1.036532863532586
$ws.Range("E25").Value = 1.04225839362846
$ws.Range("F25").Value = 1.050124374715419
$ws.Range("I25").Value = 1.036456374819708
$ws.Range("J25").Value = 1.039921237571944
$ws.Range("K25").Value = 1.039583091414209
$ws.Range("L25").Value = 1.045290763391571
$ws.Range("M25").Value = 1.053132556724184
